# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Price cells whose new values look numeric,
# so Excel keeps them stored as text (preserving e.g. trailing zeros)
# instead of silently converting the assigned string into a number.
# (Applied per-cell individually since multi-area comma ranges only
# reliably affect the first area in this COM runtime.)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values cell by cell
$ws.Range("D2").Value = '62.803.64'
$ws.Range("E2").Value = '  -1.30%  '
$ws.Range("D3").Value = '2.677.36'
$ws.Range("E3").Value = '  -1.83%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '553.52'
$ws.Range("E5").Value = '  -1.84%  '
$ws.Range("D6").Value = '157.04'
$ws.Range("E6").Value = '  -1.79%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '0.588'
$ws.Range("E8").Value = '  -1.39%  '
$ws.Range("E9").Value = '  -3.77%  '
$ws.Range("D10").Value = '0.161'
$ws.Range("E10").Value = '  -3.10%  '
$ws.Range("E11").Value = '  -2.41%  '
$ws.Range("D12").Value = '0.365'
$ws.Range("E12").Value = '  -3.75%  '
$ws.Range("D13").Value = '3.152.52'
$ws.Range("E13").Value = '  -1.76%  '
$ws.Range("D14").Value = '26.11'
$ws.Range("E14").Value = '  -2.53%  '
$ws.Range("D15").Value = '62.713.92'
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("E16").Value = '  -2.74%  '
$ws.Range("D17").Value = '2.681.94'
$ws.Range("E17").Value = '  -1.76%  '
$ws.Range("D18").Value = '11.80'
$ws.Range("E18").Value = '  -6.19%  '
$ws.Range("D19").Value = '4.58'
$ws.Range("E19").Value = '  -3.42%  '
$ws.Range("D20").Value = '343.73'
$ws.Range("E20").Value = '  -2.91%  '
$ws.Range("D21").Value = '6.19'
$ws.Range("E21").Value = '  -5.51%  '
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").Value = '0.508'
$ws.Range("E23").Value = '  -2.59%  '
$ws.Range("D24").Value = '63.23'
$ws.Range("E24").Value = '  -1.84%  '
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").Value = '8.13'
$ws.Range("E27").Value = '  -2.94%  '
$ws.Range("D28").Value = '1.40'
$ws.Range("E28").Value = '  +4.66%  '
$ws.Range("D29").Value = '0.0₃0848'
$ws.Range("E29").Value = '  -6.21%  '
$ws.Range("D30").Value = '7.25'
$ws.Range("E30").Value = '  +0.84%  '
$ws.Range("D31").Value = '1.92'
$ws.Range("E31").Value = '  -2.26%  '
$ws.Range("D32").Value = '162.78'
$ws.Range("E32").Value = '  -2.11%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = '4.84'
$ws.Range("E34").Value = '  -1.34%  '
$ws.Range("D35").Value = '1.45'
$ws.Range("E35").Value = '  -1.86%  '
$ws.Range("D36").Value = '19.42'
$ws.Range("E36").Value = '  -3.11%  '
$ws.Range("D37").Value = '1.78'
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("D38").Value = '337.80'
$ws.Range("E38").Value = '  -2.06%  '
$ws.Range("D39").Value = '6.12'
$ws.Range("E39").Value = '  -2.70%  '
$ws.Range("D40").Value = '0.922'
$ws.Range("E40").Value = '  -5.51%  '
$ws.Range("D41").Value = '3.98'
$ws.Range("E41").Value = '  -2.66%  '
$ws.Range("D42").Value = '38.32'
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("D43").Value = '20.77'
$ws.Range("E43").Value = '  -4.90%  '
$ws.Range("D44").Value = '20.12'
$ws.Range("E44").Value = '  -4.66%  '
$ws.Range("E45").Value = '  -1.64%  '
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").Value = '0.0555'
$ws.Range("E47").Value = '  -4.82%  '
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").Value = '11.00'
$ws.Range("E48").Value = '  -0.68%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '130.04'
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("D50").Value = '0.0968'
$ws.Range("E50").Value = '  -3.08%  '
$ws.Range("D51").Value = '0.0239'
$ws.Range("E51").Value = '  -4.34%  '
